$wb = $excel.ActiveWorkbook

# Sheet: ALC (sheet index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 8).Value = 96.875  # H4: 117.666664 -> 96.875
$ws.Cells.Item(4, 9).Value = 96.875  # I4: 82.90909000000001 -> 96.875
$ws.Cells.Item(4, 10).Value = 0  # J4: 500 -> 0
$ws.Cells.Item(4, 11).Value = 96.875  # K4: 82.90909000000001 -> 96.875
$ws.Cells.Item(4, 12).Value = 0  # L4: 500 -> 0
$ws.Cells.Item(4, 13).Value = 17.125  # M4: 31.09090999999999 -> 17.125
$ws.Cells.Item(4, 14).ClearContents()  # N4: -728 -> (cleared)
$ws.Cells.Item(34, 8).Value = 4998  # H34: 6249.5 -> 4998
$ws.Cells.Item(34, 9).Value = 4998  # I34: 6249.5 -> 4998
$ws.Cells.Item(34, 11).Value = 4998  # K34: 6249.5 -> 4998
$ws.Cells.Item(34, 13).Value = -4795  # M34: -6046.5 -> -4795
$ws.Cells.Item(36, 8).Value = 4998  # H36: 6249.5 -> 4998
$ws.Cells.Item(36, 9).Value = 4998  # I36: 6249.5 -> 4998
$ws.Cells.Item(36, 11).Value = 4998  # K36: 6249.5 -> 4998
$ws.Cells.Item(36, 13).Value = -4283  # M36: -5534.5 -> -4283
$ws.Cells.Item(99, 8).Value = 277.42856  # H99: 304.16666 -> 277.42856
$ws.Cells.Item(99, 9).Value = 290.5  # I99: 308.8 -> 290.5
$ws.Cells.Item(99, 10).Value = 199  # J99: 281 -> 199
$ws.Cells.Item(99, 11).Value = 871.5  # K99: 926.4000000000001 -> 871.5
$ws.Cells.Item(99, 12).Value = 597  # L99: 843 -> 597
$ws.Cells.Item(99, 13).Value = 626.5  # M99: 571.5999999999999 -> 626.5
$ws.Cells.Item(99, 14).Value = -3593  # N99: -3839 -> -3593
$ws.Cells.Item(121, 8).Value = 1871.625  # H121: 1977.0667 -> 1871.625
$ws.Cells.Item(121, 10).Value = 1871.625  # J121: 1977.0667 -> 1871.625
$ws.Cells.Item(121, 12).Value = 5614.875  # L121: 5931.2001 -> 5614.875
$ws.Cells.Item(121, 14).Value = -9108.875  # N121: -9425.2001 -> -9108.875
$ws.Cells.Item(138, 8).Value = 5569031  # H138: 5197987 -> 5569031
$ws.Cells.Item(138, 10).Value = 7795353.5  # J138: 7086991.5 -> 7795353.5
$ws.Cells.Item(138, 12).Value = 23386060.5  # L138: 21260974.5 -> 23386060.5
$ws.Cells.Item(138, 14).Value = -23396340.5  # N138: -21271254.5 -> -23396340.5

# Sheet: ARM (sheet index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 26429.727  # H32: 20057.621 -> 26429.727
$ws.Cells.Item(32, 9).Value = 30405.055  # I32: 22067.078 -> 30405.055
$ws.Cells.Item(32, 11).Value = 30405.055  # K32: 22067.078 -> 30405.055
$ws.Cells.Item(32, 13).Value = -30118.055  # M32: -21780.078 -> -30118.055
$ws.Cells.Item(46, 8).Value = 5968.75  # H46: 6257.143 -> 5968.75
$ws.Cells.Item(46, 9).Value = 5749.6  # I46: 6199.5 -> 5749.6
$ws.Cells.Item(46, 11).Value = 5749.6  # K46: 6199.5 -> 5749.6
$ws.Cells.Item(46, 13).Value = -5430.6  # M46: -5880.5 -> -5430.6
$ws.Cells.Item(102, 8).Value = 1114.6666  # H102: 1580.5454 -> 1114.6666
$ws.Cells.Item(102, 9).Value = 985.0769  # I102: 1496.8889 -> 985.0769
$ws.Cells.Item(102, 11).Value = 985.0769  # K102: 1496.8889 -> 985.0769
$ws.Cells.Item(102, 13).Value = 636.9231  # M102: 125.1111000000001 -> 636.9231

# Sheet: BSM (sheet index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(107, 8).Value = 0  # H107: 1229.8572 -> 0
$ws.Cells.Item(107, 9).Value = 0  # I107: 1229.8572 -> 0
$ws.Cells.Item(107, 11).Value = 0  # K107: 1229.8572 -> 0
$ws.Cells.Item(107, 13).ClearContents()  # M107: 690.1428000000001 -> (cleared)

# Sheet: CRP (sheet index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(86, 8).Value = 18540.889  # H86: 16806.092 -> 18540.889
$ws.Cells.Item(86, 9).Value = 22221.75  # I86: 17814.334 -> 22221.75
$ws.Cells.Item(86, 11).Value = 22221.75  # K86: 17814.334 -> 22221.75
$ws.Cells.Item(86, 13).Value = -21098.75  # M86: -16691.334 -> -21098.75
$ws.Cells.Item(89, 8).Value = 18540.889  # H89: 16806.092 -> 18540.889
$ws.Cells.Item(89, 9).Value = 22221.75  # I89: 17814.334 -> 22221.75
$ws.Cells.Item(89, 11).Value = 111108.75  # K89: 89071.67 -> 111108.75
$ws.Cells.Item(89, 13).Value = -105492.75  # M89: -83455.67 -> -105492.75
$ws.Cells.Item(99, 8).Value = 4273.75  # H99: 3491.6428 -> 4273.75
$ws.Cells.Item(99, 9).Value = 4299.3335  # I99: 3370.1428 -> 4299.3335
$ws.Cells.Item(99, 10).Value = 4258.4  # J99: 3613.1428 -> 4258.4
$ws.Cells.Item(99, 11).Value = 4299.3335  # K99: 3370.1428 -> 4299.3335
$ws.Cells.Item(99, 12).Value = 4258.4  # L99: 3613.1428 -> 4258.4
$ws.Cells.Item(99, 13).Value = -2801.3335  # M99: -1872.1428 -> -2801.3335
$ws.Cells.Item(99, 14).Value = -7254.4  # N99: -6609.1428 -> -7254.4
$ws.Cells.Item(105, 8).Value = 5156.857  # H105: 3518.9092 -> 5156.857
$ws.Cells.Item(105, 9).Value = 5366  # I105: 3285.2856 -> 5366
$ws.Cells.Item(105, 10).Value = 5000  # J105: 3927.75 -> 5000
$ws.Cells.Item(105, 11).Value = 5366  # K105: 3285.2856 -> 5366
$ws.Cells.Item(105, 12).Value = 5000  # L105: 3927.75 -> 5000
$ws.Cells.Item(105, 13).Value = -3619  # M105: -1538.2856 -> -3619
$ws.Cells.Item(105, 14).Value = -8494  # N105: -7421.75 -> -8494
$ws.Cells.Item(107, 8).Value = 593.625  # H107: 619.2 -> 593.625
$ws.Cells.Item(107, 9).Value = 382.81818  # I107: 400.1 -> 382.81818
$ws.Cells.Item(107, 11).Value = 382.81818  # K107: 400.1 -> 382.81818
$ws.Cells.Item(107, 13).Value = 1537.18182  # M107: 1519.9 -> 1537.18182
$ws.Cells.Item(122, 8).Value = 0  # H122: 1001.0526 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 833.3333 -> 0
$ws.Cells.Item(122, 10).Value = 0  # J122: 1288.5714 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 2499.9999 -> 0
$ws.Cells.Item(122, 12).Value = 0  # L122: 3865.7142 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: -49.9998999999998 -> (cleared)
$ws.Cells.Item(122, 14).ClearContents()  # N122: -8765.7142 -> (cleared)
$ws.Cells.Item(126, 8).Value = 4273.75  # H126: 3491.6428 -> 4273.75
$ws.Cells.Item(126, 9).Value = 4299.3335  # I126: 3370.1428 -> 4299.3335
$ws.Cells.Item(126, 10).Value = 4258.4  # J126: 3613.1428 -> 4258.4
$ws.Cells.Item(126, 11).Value = 12898.0005  # K126: 10110.4284 -> 12898.0005
$ws.Cells.Item(126, 12).Value = 12775.2  # L126: 10839.4284 -> 12775.2
$ws.Cells.Item(126, 13).Value = -10428.0005  # M126: -7640.428400000001 -> -10428.0005
$ws.Cells.Item(126, 14).Value = -17715.2  # N126: -15779.4284 -> -17715.2
$ws.Cells.Item(132, 8).Value = 1410.4375  # H132: 1514.9231 -> 1410.4375
$ws.Cells.Item(132, 9).Value = 1410.4375  # I132: 1514.9231 -> 1410.4375
$ws.Cells.Item(132, 11).Value = 4231.3125  # K132: 4544.7693 -> 4231.3125
$ws.Cells.Item(132, 13).Value = -1701.3125  # M132: -2014.7693 -> -1701.3125

# Sheet: CUL (sheet index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 763.26086  # H5: 768.6818 -> 763.26086
$ws.Cells.Item(5, 10).Value = 801.3570999999999  # J5: 813.46155 -> 801.3570999999999
$ws.Cells.Item(5, 12).Value = 2404.0713  # L5: 2440.38465 -> 2404.0713
$ws.Cells.Item(5, 14).Value = -2628.0713  # N5: -2664.38465 -> -2628.0713
$ws.Cells.Item(68, 8).Value = 3049.0728  # H68: 3054.182 -> 3049.0728
$ws.Cells.Item(68, 10).Value = 3245.3264  # J68: 3251.0613 -> 3245.3264
$ws.Cells.Item(68, 12).Value = 9735.9792  # L68: 9753.1839 -> 9735.9792
$ws.Cells.Item(68, 14).Value = -11357.9792  # N68: -11375.1839 -> -11357.9792
$ws.Cells.Item(71, 8).Value = 3049.0728  # H71: 3054.182 -> 3049.0728
$ws.Cells.Item(71, 10).Value = 3245.3264  # J71: 3251.0613 -> 3245.3264
$ws.Cells.Item(71, 12).Value = 29207.9376  # L71: 29259.5517 -> 29207.9376
$ws.Cells.Item(71, 14).Value = -37319.9376  # N71: -37371.5517 -> -37319.9376
$ws.Cells.Item(135, 8).Value = 763.26086  # H135: 768.6818 -> 763.26086
$ws.Cells.Item(135, 10).Value = 801.3570999999999  # J135: 813.46155 -> 801.3570999999999
$ws.Cells.Item(135, 12).Value = 7212.2139  # L135: 7321.15395 -> 7212.2139
$ws.Cells.Item(135, 14).Value = -12282.2139  # N135: -12391.15395 -> -12282.2139

# Sheet: GSM (sheet index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 1554.2273  # H102: 1683.3684 -> 1554.2273
$ws.Cells.Item(102, 9).Value = 1326.1052  # I102: 1469.7059 -> 1326.1052
$ws.Cells.Item(102, 10).Value = 2999  # J102: 3499.5 -> 2999
$ws.Cells.Item(102, 11).Value = 1326.1052  # K102: 1469.7059 -> 1326.1052
$ws.Cells.Item(102, 12).Value = 2999  # L102: 3499.5 -> 2999
$ws.Cells.Item(102, 13).Value = 295.8948  # M102: 152.2941000000001 -> 295.8948
$ws.Cells.Item(102, 14).Value = -6243  # N102: -6743.5 -> -6243
$ws.Cells.Item(113, 8).Value = 2102.4546  # H113: 2430.3635 -> 2102.4546
$ws.Cells.Item(113, 9).Value = 2123.7778  # I113: 2192.7778 -> 2123.7778
$ws.Cells.Item(113, 10).Value = 2006.5  # J113: 3499.5 -> 2006.5
$ws.Cells.Item(113, 11).Value = 2123.7778  # K113: 2192.7778 -> 2123.7778
$ws.Cells.Item(113, 12).Value = 2006.5  # L113: 3499.5 -> 2006.5
$ws.Cells.Item(113, 13).Value = 46.22220000000016  # M113: -22.77779999999984 -> 46.22220000000016
$ws.Cells.Item(113, 14).Value = -6346.5  # N113: -7839.5 -> -6346.5
$ws.Cells.Item(132, 8).Value = 6688.9165  # H132: 5920 -> 6688.9165
$ws.Cells.Item(132, 9).Value = 7848.8125  # I132: 6841.737 -> 7848.8125
$ws.Cells.Item(132, 10).Value = 4369.125  # J132: 3974.111 -> 4369.125
$ws.Cells.Item(132, 11).Value = 23546.4375  # K132: 20525.211 -> 23546.4375
$ws.Cells.Item(132, 12).Value = 13107.375  # L132: 11922.333 -> 13107.375
$ws.Cells.Item(132, 13).Value = -21016.4375  # M132: -17995.211 -> -21016.4375
$ws.Cells.Item(132, 14).Value = -18167.375  # N132: -16982.333 -> -18167.375
$ws.Cells.Item(134, 8).Value = 95162.5  # H134: 99999 -> 95162.5
$ws.Cells.Item(134, 10).Value = 95162.5  # J134: 99999 -> 95162.5
$ws.Cells.Item(134, 12).Value = 285487.5  # L134: 299997 -> 285487.5
$ws.Cells.Item(134, 14).Value = -290557.5  # N134: -305067 -> -290557.5

# Sheet: LTW (sheet index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 2557.7715  # H40: 2703.875 -> 2557.7715
$ws.Cells.Item(40, 9).Value = 2440.0605  # I40: 2584.1333 -> 2440.0605
$ws.Cells.Item(40, 11).Value = 2440.0605  # K40: 2584.1333 -> 2440.0605
$ws.Cells.Item(40, 13).Value = -2304.0605  # M40: -2448.1333 -> -2304.0605
$ws.Cells.Item(46, 8).Value = 5043  # H46: 4515.4736 -> 5043
$ws.Cells.Item(46, 9).Value = 0  # I46: 999 -> 0
$ws.Cells.Item(46, 10).Value = 5043  # J46: 4710.8335 -> 5043
$ws.Cells.Item(46, 11).Value = 0  # K46: 999 -> 0
$ws.Cells.Item(46, 12).Value = 5043  # L46: 4710.8335 -> 5043
$ws.Cells.Item(46, 13).ClearContents()  # M46: -811 -> (cleared)
$ws.Cells.Item(46, 14).Value = -5419  # N46: -5086.8335 -> -5419
$ws.Cells.Item(68, 8).Value = 2216.6667  # H68: 2118.182 -> 2216.6667
$ws.Cells.Item(68, 9).Value = 2216.6667  # I68: 2125 -> 2216.6667
$ws.Cells.Item(68, 10).Value = 0  # J68: 2100 -> 0
$ws.Cells.Item(68, 11).Value = 2216.6667  # K68: 2125 -> 2216.6667
$ws.Cells.Item(68, 12).Value = 0  # L68: 2100 -> 0
$ws.Cells.Item(68, 13).Value = -1467.6667  # M68: -1376 -> -1467.6667
$ws.Cells.Item(68, 14).ClearContents()  # N68: -3598 -> (cleared)
$ws.Cells.Item(71, 8).Value = 2216.6667  # H71: 2118.182 -> 2216.6667
$ws.Cells.Item(71, 9).Value = 2216.6667  # I71: 2125 -> 2216.6667
$ws.Cells.Item(71, 10).Value = 0  # J71: 2100 -> 0
$ws.Cells.Item(71, 11).Value = 11083.3335  # K71: 10625 -> 11083.3335
$ws.Cells.Item(71, 12).Value = 0  # L71: 10500 -> 0
$ws.Cells.Item(71, 13).Value = -7339.333500000001  # M71: -6881 -> -7339.333500000001
$ws.Cells.Item(71, 14).ClearContents()  # N71: -17988 -> (cleared)
$ws.Cells.Item(82, 8).Value = 3405.55  # H82: 4343.3335 -> 3405.55
$ws.Cells.Item(82, 9).Value = 2703.6924  # I82: 4172.846 -> 2703.6924
$ws.Cells.Item(82, 10).Value = 4709  # J82: 4620.375 -> 4709
$ws.Cells.Item(82, 11).Value = 2703.6924  # K82: 4172.846 -> 2703.6924
$ws.Cells.Item(82, 12).Value = 4709  # L82: 4620.375 -> 4709
$ws.Cells.Item(82, 13).Value = -2342.6924  # M82: -3811.846 -> -2342.6924
$ws.Cells.Item(82, 14).Value = -5431  # N82: -5342.375 -> -5431
$ws.Cells.Item(85, 8).Value = 3405.55  # H85: 4343.3335 -> 3405.55
$ws.Cells.Item(85, 9).Value = 2703.6924  # I85: 4172.846 -> 2703.6924
$ws.Cells.Item(85, 10).Value = 4709  # J85: 4620.375 -> 4709
$ws.Cells.Item(85, 11).Value = 2703.6924  # K85: 4172.846 -> 2703.6924
$ws.Cells.Item(85, 12).Value = 4709  # L85: 4620.375 -> 4709
$ws.Cells.Item(85, 13).Value = -1455.6924  # M85: -2924.846 -> -1455.6924
$ws.Cells.Item(85, 14).Value = -7205  # N85: -7116.375 -> -7205
$ws.Cells.Item(122, 8).Value = 2149.9285  # H122: 2391 -> 2149.9285
$ws.Cells.Item(122, 9).Value = 1980  # I122: 2162.25 -> 1980
$ws.Cells.Item(122, 10).Value = 2574.75  # J122: 2848.5 -> 2574.75
$ws.Cells.Item(122, 11).Value = 5940  # K122: 6486.75 -> 5940
$ws.Cells.Item(122, 12).Value = 7724.25  # L122: 8545.5 -> 7724.25
$ws.Cells.Item(122, 13).Value = -3490  # M122: -4036.75 -> -3490
$ws.Cells.Item(122, 14).Value = -12624.25  # N122: -13445.5 -> -12624.25
$ws.Cells.Item(132, 8).Value = 3806.756  # H132: 3802.2 -> 3806.756
$ws.Cells.Item(132, 10).Value = 5019.4614  # J132: 5105.3335 -> 5019.4614
$ws.Cells.Item(132, 12).Value = 15058.3842  # L132: 15316.0005 -> 15058.3842
$ws.Cells.Item(132, 14).Value = -20118.3842  # N132: -20376.0005 -> -20118.3842

# Sheet: WVR (sheet index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(14, 8).Value = 0  # H14: 10000 -> 0
$ws.Cells.Item(14, 9).Value = 0  # I14: 10000 -> 0
$ws.Cells.Item(14, 11).Value = 0  # K14: 10000 -> 0
$ws.Cells.Item(14, 13).ClearContents()  # M14: -9832 -> (cleared)
$ws.Cells.Item(62, 8).Value = 10990.462  # H62: 11741 -> 10990.462
$ws.Cells.Item(62, 9).Value = 11055.429  # I62: 12065.667 -> 11055.429
$ws.Cells.Item(62, 10).Value = 10914.667  # J62: 11497.5 -> 10914.667
$ws.Cells.Item(62, 11).Value = 11055.429  # K62: 12065.667 -> 11055.429
$ws.Cells.Item(62, 12).Value = 10914.667  # L62: 11497.5 -> 10914.667
$ws.Cells.Item(62, 13).Value = -10431.429  # M62: -11441.667 -> -10431.429
$ws.Cells.Item(62, 14).Value = -12162.667  # N62: -12745.5 -> -12162.667
$ws.Cells.Item(65, 8).Value = 10990.462  # H65: 11741 -> 10990.462
$ws.Cells.Item(65, 9).Value = 11055.429  # I65: 12065.667 -> 11055.429
$ws.Cells.Item(65, 10).Value = 10914.667  # J65: 11497.5 -> 10914.667
$ws.Cells.Item(65, 11).Value = 55277.145  # K65: 60328.335 -> 55277.145
$ws.Cells.Item(65, 12).Value = 54573.335  # L65: 57487.5 -> 54573.335
$ws.Cells.Item(65, 13).Value = -52157.145  # M65: -57208.335 -> -52157.145
$ws.Cells.Item(65, 14).Value = -60813.335  # N65: -63727.5 -> -60813.335
$ws.Cells.Item(122, 8).Value = 14289230  # H122: 11366601 -> 14289230
$ws.Cells.Item(122, 9).Value = 20003666  # I122: 14708796 -> 20003666
$ws.Cells.Item(122, 11).Value = 60010998  # K122: 44126388 -> 60010998
$ws.Cells.Item(122, 13).Value = -60008548  # M122: -44123938 -> -60008548
$ws.Cells.Item(126, 8).Value = 2394.7932  # H126: 2395 -> 2394.7932
$ws.Cells.Item(126, 9).Value = 1846.4546  # I126: 1808.4783 -> 1846.4546
$ws.Cells.Item(126, 10).Value = 4118.143  # J126: 4643.3335 -> 4118.143
$ws.Cells.Item(126, 11).Value = 5539.3638  # K126: 5425.4349 -> 5539.3638
$ws.Cells.Item(126, 12).Value = 12354.429  # L126: 13930.0005 -> 12354.429
$ws.Cells.Item(126, 13).Value = -3069.3638  # M126: -2955.4349 -> -3069.3638
$ws.Cells.Item(126, 14).Value = -17294.429  # N126: -18870.0005 -> -17294.429
$ws.Cells.Item(132, 8).Value = 11908768  # H132: 12504207 -> 11908768
$ws.Cells.Item(132, 9).Value = 11908768  # I132: 12504207 -> 11908768
$ws.Cells.Item(132, 11).Value = 35726304  # K132: 37512621 -> 35726304
$ws.Cells.Item(132, 13).Value = -35723774  # M132: -37510091 -> -35723774
$ws.Cells.Item(136, 8).Value = 5255  # H136: 5353.46 -> 5255
$ws.Cells.Item(136, 10).Value = 3062.7778  # J136: 3223.4119 -> 3062.7778
$ws.Cells.Item(136, 12).Value = 9188.3334  # L136: 9670.235700000001 -> 9188.3334
$ws.Cells.Item(136, 14).Value = -14288.3334  # N136: -14770.2357 -> -14288.3334
